$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.726.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "'2.473.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'586.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").Value = "'175.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.85%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.53%  "

$ws.Range("D9").Value = "'0.143"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.55%  "

$ws.Range("E10").Value = "  -1.63%  "

$ws.Range("D11").Value = "'4.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").Value = "'0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("D13").Value = "'2.926.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "'25.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "'67.807.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("D16").Value = "'0.0000169"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").Value = "'2.494.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "'10.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.03%  "

$ws.Range("D19").Value = "'7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.22%  "

$ws.Range("D20").Value = "'347.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.18%  "

$ws.Range("D21").Value = "'4.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "'70.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("D24").Value = "'4.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.67%  "

$ws.Range("D25").Value = "'1.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.71%  "

$ws.Range("D26").Value = "'8.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.10%  "

$ws.Range("D27").Value = "'2.584.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").Value = "'0.0₃0887"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.46%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'7.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'492.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.94%  "

$ws.Range("D32").Value = "'1.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "

$ws.Range("D33").Value = "'1.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "'163.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").Value = "'0.119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").Value = "'18.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "'1.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.88%  "

$ws.Range("D41").Value = "'1.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "

$ws.Range("D42").Value = "'0.325"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("D43").Value = "'4.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "

$ws.Range("D44").Value = "'2.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.53%  "

$ws.Range("D45").Value = "'146.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.38%  "

$ws.Range("D46").Value = "'3.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.67%  "

$ws.Range("D47").Value = "'0.508"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").Value = "'0.0₆0251"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.76%  "

$ws.Range("D49").Value = "'0.0735"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("D50").Value = "'1.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").Value = "'0.574"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.83%  "
